$d = $word.ActiveDocument

# --- 1) Clear the student RA number in the body ------------------------
# " 000110343434 - 2 "  ->  "  " (two spaces)
$d.Content.Find.Execute(" 000110343434 - 2 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "  ", 2)

# --- 2) Replace the "QWR" placeholder in the body ("A QWR," sentence) --
$d.Content.Find.Execute("QWR", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TERE", 2)

# --- 3) Header text replacements ----------------------------------------
# Grab the primary header range for the (single) section.
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdrRange = $hdr.Range

# "QWER" -> "TRE"  (DIRETORIA DE ENSINO REGIAO ...)
$hdrRange.Find.Execute("QWER", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "TRE", 2)

# "QWR" -> "TERE"  (... - DEP.)
$hdrRange.Find.Execute("QWR", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "TERE", 2)

# "Qwer" -> "Tre"  (address line, 5 occurrences)
$hdrRange.Find.Execute("Qwer", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "Tre", 2)

# "qwer" -> "tre"  (CEP / Tel / Email lines, 3 occurrences)
$hdrRange.Find.Execute("qwer", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "tre", 2)
